$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: quality_comparison
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

# C1 gets a top+bottom border (no left/right) - new style referencing borderId 4
$s1c1 = $ws1.Range("C1")
$s1c1.ClearFormats()
$s1c1.Borders.LineStyle = 1              # all four sides thin
$s1c1.Borders.Item(7).LineStyle = -4142  # xlEdgeLeft -> none
$s1c1.Borders.Item(10).LineStyle = -4142 # xlEdgeRight -> none

# D1 gets a top+right+bottom border (no left) - new style referencing borderId 5
$s1d1 = $ws1.Range("D1")
$s1d1.ClearFormats()
$s1d1.Borders.LineStyle = 1              # all four sides thin
$s1d1.Borders.Item(7).LineStyle = -4142  # xlEdgeLeft -> none

# C2 text: fedcore -> approach
$ws1.Range("C2").Value = "approach"

# ---------------------------------------------------------------------------
# Sheet 2: computational_comparison
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

# C1 / F1 mirror sheet1's C1 formatting (top+bottom border, no left/right)
$s1c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

# D1 / G1 mirror sheet1's D1 formatting (top+right+bottom border, no left)
$s1d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# C2 / F2 text: fedcore -> approach
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()

Write-Host "done"
